# Natmi following Dr Hou advice
# Update the LR-pair stats for Clec11a-Itga10 after recomputation:
# ligand- and receptor-expressing cell counts go from 1 to 3, which in turn
# changes the dependent average/total expression values and derived
# specificity scores for every data row (rows 2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{ E=3; G=9.930652333333333;  H=29.791957;        I=0.9673539331442913; J=0.9673539331442912; K=3; M=0.4043423333333334; N=1.213027;         O=0.08238992852068679; P=0.08238992852068679; Q=4.015383135982111;  R=36.138448223839;    S=0.07970022140596339;  T=0.07970022140596338 }
    3 = @{ E=3; G=9.930652333333333;  H=29.791957;        I=0.9673539331442913; J=0.9673539331442912; K=3; M=1.358024333333333;  N=4.074073;         O=0.276714849099039;   P=0.276714849099039;   Q=13.48606751454011;  R=121.374607630861;   S=0.2676811976353844;   T=0.2676811976353844 }
    4 = @{ E=3; G=9.930652333333333;  H=29.791957;        I=0.9673539331442913; J=0.9673539331442912; K=3; M=3.145300333333333;  N=9.435900999999999; O=0.6408952223802742;  P=0.6408952223802742;  Q=31.23488409425077;  R=281.113956848257;   S=0.6199725141029435;   T=0.6199725141029434 }
    5 = @{ E=3; G=0.3351376666666666; H=1.005413;         I=0.03264606685570878; J=0.03264606685570878; K=3; M=0.4043423333333334; N=1.213027;         O=0.08238992852068679; P=0.08238992852068679; Q=0.1355103461278889; R=1.219593115151;     S=0.002689707114723409; T=0.002689707114723408 }
    6 = @{ E=3; G=0.3351376666666666; H=1.005413;         I=0.03264606685570878; J=0.03264606685570878; K=3; M=1.358024333333333;  N=4.074073;         O=0.276714849099039;   P=0.276714849099039;   Q=0.4551251063498889; R=4.096125957149;     S=0.009033651463654596; T=0.009033651463654592 }
    7 = @{ E=3; G=0.3351376666666666; H=1.005413;         I=0.03264606685570878; J=0.03264606685570878; K=3; M=3.145300333333333;  N=9.435900999999999; O=0.6408952223802742;  P=0.6408952223802742;  Q=1.054108614679222;  R=9.486977532112999;  S=0.02092270827733078;  T=0.02092270827733077 }
}

foreach ($row in $rowData.Keys) {
    $cols = $rowData[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
